# Edit script: update transmittal date and rewrite the main body paragraph
# to reflect that the case file is now escorted by police (instead of a
# complaint submission), per the commit's updated wording.

$d = $word.ActiveDocument

# --- Change the transmittal date "04/10/2025" -> "03/11/2025".
#     It occurs twice (top-right header block and bottom-right signature
#     block); wdReplaceAll (2) updates every occurrence in one pass. ---
$d.Content.Find.Execute("04/10/2025", $true, $false, $false, $false, $false, $true, 1, $false, "03/11/2025", 2) | Out-Null

# --- Rewrite the main body paragraph (the one starting with
#     "Υποβάλλεται συννημένα ...") with the new wording. The old paragraph
#     contained a line break + tab followed by a second sentence; the new
#     text is a single sentence, so we replace the whole paragraph content
#     (keeping the trailing paragraph mark intact). ---
$newBodyText = '   Αποστέλλεται συνοδεία Αστυνομικών Υπηρεσίας μας και με την σε βάρος του σχηματισθείσα Δικογραφία ο {{surnamePerperator}} {{namePerperator}} του {{fathernamePerperator}} και της {{mothernamePerperator}}γεν. {{dateOfBirthPerperator }} στη {{placeOfBirthPerperator }} κατ.{{ addressPerperator }},αριθμός τηλεφώνου {{ telPreperator }}, ηλεκτρονικό ταχυδρομείου{{ emailPreperator }}, κάτοχος του υπ αριθμόν {{ DATperperator}} που εκδόθηκε την {{ issuedPerperator }} από {{place_issuedPerperator}} Α.Φ.Μ : {{afmPreperator}}, Δ.Ο.Υ : {{ doyPrep }}, κατηγορούμενος για παραβάσεις του/τωνάρθρων {{ offences}} πράξεις που έλαβαν χώρα  στη {{placeOfCrime}} στις {{ dateOfCrime}} και περί ώρα {{ hourOfCrime }}'

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.IndexOf("Υποβάλλεται συννημένα") -ge 0) {
        $target = $p
        break
    }
}
if ($null -eq $target) {
    throw "Could not locate the body paragraph to replace"
}

$bodyRange = $target.Range
$bodyRange.MoveEnd(1, -1) | Out-Null   # exclude the trailing paragraph mark
$bodyRange.Text = $newBodyText
